$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Fill in row 10 with the new test-data entry
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Delete For Agency Object Profile"
$ws.Range("C10").Value = "AutoTitleEdited"

# Update the active cell selection to C13
$ws.Range("C13").Select()
